$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3280
$ws.Cells.Item(40, 9).Value = 1900
$ws.Cells.Item(40, 10).Value = 3625
$ws.Cells.Item(40, 11).Value = 1900
$ws.Cells.Item(40, 12).Value = 3625
$ws.Cells.Item(40, 13).Value = -1725
$ws.Cells.Item(40, 14).Value = -3975
$ws.Cells.Item(42, 8).Value = 478.33334
$ws.Cells.Item(42, 9).Value = 413.125
$ws.Cells.Item(42, 11).Value = 1239.375
$ws.Cells.Item(42, 13).Value = -1009.375
$ws.Cells.Item(86, 8).Value = 2164.6
$ws.Cells.Item(86, 9).Value = 1803.1428
$ws.Cells.Item(86, 10).Value = 3008
$ws.Cells.Item(86, 11).Value = 1803.1428
$ws.Cells.Item(86, 12).Value = 3008
$ws.Cells.Item(86, 13).Value = -680.1428000000001
$ws.Cells.Item(86, 14).Value = -5254
$ws.Cells.Item(89, 8).Value = 2164.6
$ws.Cells.Item(89, 9).Value = 1803.1428
$ws.Cells.Item(89, 10).Value = 3008
$ws.Cells.Item(89, 11).Value = 9015.714
$ws.Cells.Item(89, 12).Value = 15040
$ws.Cells.Item(89, 13).Value = -3399.714
$ws.Cells.Item(89, 14).Value = -26272
$ws.Cells.Item(92, 8).Value = 2479.4707
$ws.Cells.Item(92, 9).Value = 186.16667
$ws.Cells.Item(92, 10).Value = 7983.4
$ws.Cells.Item(92, 11).Value = 186.16667
$ws.Cells.Item(92, 12).Value = 7983.4
$ws.Cells.Item(92, 13).Value = 1061.83333
$ws.Cells.Item(92, 14).Value = -10479.4
$ws.Cells.Item(100, 8).Value = 6552.4165
$ws.Cells.Item(100, 9).Value = 6806.8184
$ws.Cells.Item(100, 11).Value = 6806.8184
$ws.Cells.Item(100, 13).Value = -6265.8184
$ws.Cells.Item(113, 8).Value = 9007.615
$ws.Cells.Item(113, 9).Value = 8637.375
$ws.Cells.Item(113, 11).Value = 8637.375
$ws.Cells.Item(113, 13).Value = -5383.375
$ws.Cells.Item(137, 8).Value = 31728.9
$ws.Cells.Item(137, 9).Value = 75749.5
$ws.Cells.Item(137, 10).Value = 2381.8333
$ws.Cells.Item(137, 11).Value = 227248.5
$ws.Cells.Item(137, 12).Value = 7145.499899999999
$ws.Cells.Item(137, 13).Value = -224698.5
$ws.Cells.Item(137, 14).Value = -12245.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 19759.871
$ws.Cells.Item(32, 9).Value = 20834.02
$ws.Cells.Item(32, 11).Value = 20834.02
$ws.Cells.Item(32, 13).Value = -20547.02
$ws.Cells.Item(45, 8).Value = 4701
$ws.Cells.Item(45, 10).Value = 5562.2
$ws.Cells.Item(45, 12).Value = 5562.2
$ws.Cells.Item(45, 14).Value = -6316.2
$ws.Cells.Item(122, 8).Value = 2914.8462
$ws.Cells.Item(122, 9).Value = 2589.4
$ws.Cells.Item(122, 11).Value = 7768.200000000001
$ws.Cells.Item(122, 13).Value = -5318.200000000001
$ws.Cells.Item(132, 8).Value = 2073.1052
$ws.Cells.Item(132, 9).Value = 1910.5
$ws.Cells.Item(132, 10).Value = 5000
$ws.Cells.Item(132, 11).Value = 5731.5
$ws.Cells.Item(132, 12).Value = 15000
$ws.Cells.Item(132, 13).Value = -3201.5
$ws.Cells.Item(132, 14).Value = -20060

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 14).Value = $null
$ws.Cells.Item(20, 8).Value = 3726.1765
$ws.Cells.Item(20, 9).Value = 4588.625
$ws.Cells.Item(20, 10).Value = 1656.3
$ws.Cells.Item(20, 11).Value = 4588.625
$ws.Cells.Item(20, 12).Value = 1656.3
$ws.Cells.Item(20, 13).Value = -4341.625
$ws.Cells.Item(20, 14).Value = -2150.3
$ws.Cells.Item(99, 8).Value = 1710.125
$ws.Cells.Item(99, 9).Value = 1383
$ws.Cells.Item(99, 10).Value = 4000
$ws.Cells.Item(99, 11).Value = 1383
$ws.Cells.Item(99, 12).Value = 4000
$ws.Cells.Item(99, 13).Value = 115
$ws.Cells.Item(99, 14).Value = -6996

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 371
$ws.Cells.Item(10, 9).Value = 371
$ws.Cells.Item(10, 11).Value = 371
$ws.Cells.Item(10, 13).Value = -232
$ws.Cells.Item(134, 8).Value = 3382.6
$ws.Cells.Item(134, 9).Value = 3139.4
$ws.Cells.Item(134, 10).Value = 3625.8
$ws.Cells.Item(134, 11).Value = 9418.200000000001
$ws.Cells.Item(134, 12).Value = 10877.4
$ws.Cells.Item(134, 13).Value = -6883.200000000001
$ws.Cells.Item(134, 14).Value = -15947.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 696.75
$ws.Cells.Item(8, 9).Value = 696.75
$ws.Cells.Item(8, 11).Value = 2090.25
$ws.Cells.Item(8, 13).Value = -1951.25
$ws.Cells.Item(98, 8).Value = 1233.3334
$ws.Cells.Item(98, 10).Value = 3000
$ws.Cells.Item(98, 12).Value = 9000
$ws.Cells.Item(98, 14).Value = -11996
$ws.Cells.Item(113, 8).Value = 869.5714
$ws.Cells.Item(113, 9).Value = 693
$ws.Cells.Item(113, 10).Value = 899
$ws.Cells.Item(113, 11).Value = 2079
$ws.Cells.Item(113, 12).Value = 2697
$ws.Cells.Item(113, 13).Value = 91
$ws.Cells.Item(113, 14).Value = -7037

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14, 8).Value = 339970
$ws.Cells.Item(14, 10).Value = 9955
$ws.Cells.Item(14, 12).Value = 9955
$ws.Cells.Item(14, 14).Value = -10291
$ws.Cells.Item(70, 8).Value = 8037.8
$ws.Cells.Item(70, 9).Value = 9173.700000000001
$ws.Cells.Item(70, 10).Value = 5766
$ws.Cells.Item(70, 11).Value = 9173.700000000001
$ws.Cells.Item(70, 12).Value = 5766
$ws.Cells.Item(70, 13).Value = -8903.700000000001
$ws.Cells.Item(70, 14).Value = -6306
$ws.Cells.Item(73, 8).Value = 8037.8
$ws.Cells.Item(73, 9).Value = 9173.700000000001
$ws.Cells.Item(73, 10).Value = 5766
$ws.Cells.Item(73, 11).Value = 9173.700000000001
$ws.Cells.Item(73, 12).Value = 5766
$ws.Cells.Item(73, 13).Value = -8237.700000000001
$ws.Cells.Item(73, 14).Value = -7638
$ws.Cells.Item(102, 8).Value = 17405.344
$ws.Cells.Item(102, 9).Value = 18749.035
$ws.Cells.Item(102, 10).Value = 4416.3335
$ws.Cells.Item(102, 11).Value = 18749.035
$ws.Cells.Item(102, 12).Value = 4416.3335
$ws.Cells.Item(102, 13).Value = -17127.035
$ws.Cells.Item(102, 14).Value = -7660.3335
$ws.Cells.Item(126, 8).Value = 1986.5834
$ws.Cells.Item(126, 9).Value = 1536.75
$ws.Cells.Item(126, 10).Value = 2886.25
$ws.Cells.Item(126, 11).Value = 4610.25
$ws.Cells.Item(126, 12).Value = 8658.75
$ws.Cells.Item(126, 13).Value = -2140.25
$ws.Cells.Item(126, 14).Value = -13598.75
$ws.Cells.Item(131, 8).Value = 21000
$ws.Cells.Item(131, 10).Value = 21000
$ws.Cells.Item(131, 12).Value = 21000
$ws.Cells.Item(131, 14).Value = -31080
$ws.Cells.Item(132, 8).Value = 2911.1365
$ws.Cells.Item(132, 9).Value = 2562.182
$ws.Cells.Item(132, 11).Value = 7686.545999999999
$ws.Cells.Item(132, 13).Value = -5156.545999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(19, 8).Value = 200
$ws.Cells.Item(19, 9).Value = 200
$ws.Cells.Item(19, 11).Value = 200
$ws.Cells.Item(19, 13).Value = -30
$ws.Cells.Item(46, 8).Value = 3398.158
$ws.Cells.Item(46, 9).Value = 1475.625
$ws.Cells.Item(46, 11).Value = 1475.625
$ws.Cells.Item(46, 13).Value = -1287.625
$ws.Cells.Item(132, 8).Value = 3866.4075
$ws.Cells.Item(132, 9).Value = 3447.476
$ws.Cells.Item(132, 10).Value = 5332.6665
$ws.Cells.Item(132, 11).Value = 10342.428
$ws.Cells.Item(132, 12).Value = 15997.9995
$ws.Cells.Item(132, 13).Value = -7812.428
$ws.Cells.Item(132, 14).Value = -21057.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(33, 8).Value = 30200
$ws.Cells.Item(33, 10).Value = 30200
$ws.Cells.Item(33, 12).Value = 30200
$ws.Cells.Item(33, 14).Value = -30700
$ws.Cells.Item(36, 8).Value = 30200
$ws.Cells.Item(36, 10).Value = 30200
$ws.Cells.Item(36, 12).Value = 30200
$ws.Cells.Item(36, 14).Value = -30700
$ws.Cells.Item(122, 8).Value = 47108.9
$ws.Cells.Item(122, 9).Value = 50204.535
$ws.Cells.Item(122, 11).Value = 150613.605
$ws.Cells.Item(122, 13).Value = -148163.605
$ws.Cells.Item(132, 8).Value = 1697.7819
$ws.Cells.Item(132, 9).Value = 1257.4419
$ws.Cells.Item(132, 10).Value = 3275.6667
$ws.Cells.Item(132, 11).Value = 3772.3257
$ws.Cells.Item(132, 12).Value = 9827.000100000001
$ws.Cells.Item(132, 13).Value = -1242.3257
$ws.Cells.Item(132, 14).Value = -14887.0001
$ws.Cells.Item(136, 8).Value = 18191
$ws.Cells.Item(136, 9).Value = 27208.428
$ws.Cells.Item(136, 11).Value = 81625.284
$ws.Cells.Item(136, 13).Value = -79075.284
